$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "cus_NPuAUX7DbheAJC"
$ws.Range("A15").Value = "cus_NPuAJsxbM30H9R"
$ws.Range("A16").Value = "cus_NPuAexM6zAUeKZ"

$ws.Range("A16").Select()
